$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.382.85'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '1.628.00'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9997'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '304.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3784'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.97'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.230'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08097'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.543'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001245'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.207'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.59%  '
$ws.Range('D17').Value = '1.628.87'
$ws.Range('E17').Value = '  -0.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06896'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.413'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').Value = '23.384.45'
$ws.Range('E23').Value = '  -0.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.255'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.447'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '149.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.280'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.316'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.48%  '
$ws.Range('D32').Value = '1.808.60'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.807'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.08'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9518'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02794'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2523'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.08843'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.103'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07171'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.356'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7056'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.30'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6449'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.325'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9991'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.997'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07992'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.201'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.67'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.32%  '
